# Consolidate the split text runs ("Below" + " " + "section-level") on
# slide 2's title into a single run by rewriting the TextRange.Text.
#
# Setting TextRange.Text to the value it already logically represents is
# treated as a no-op by the writer (no visible text change), so the runs
# would not be merged. Toggling through a different value first forces the
# text frame to be rebuilt with a single consolidated run.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)
$shape.TextFrame.TextRange.Text = " "
$shape.TextFrame.TextRange.Text = "Below section-level"
